$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.75
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73
$ws.Range("U2").Value = 1.92
$ws.Range("V2").Value = 1.77
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("V3").Value = 1.63
$ws.Range("V4").Value = 1.63
$ws.Range("I5").Value = 2.45
$ws.Range("W5").Value = 9.5
$ws.Range("AM5").Value = 26
$ws.Range("AY5").Value = 21
$ws.Range("G7").Value = 2.57
$ws.Range("I7").Value = 2.65
$ws.Range("L7").Value = 3.15
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 11
$ws.Range("P7").Value = 3.4
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.13
$ws.Range("V7").Value = 2.22
$ws.Range("W7").Value = 9.75
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 9.25
$ws.Range("Z7").Value = 30
$ws.Range("AH7").Value = 10.5
$ws.Range("AI7").Value = 16
$ws.Range("AM7").Value = 23
$ws.Range("AN7").Value = 4.65
$ws.Range("AP7").Value = 18.5
$ws.Range("AT7").Value = 2.82
$ws.Range("AU7").Value = 6.1
$ws.Range("AW7").Value = 4.75
$ws.Range("AX7").Value = 14
$ws.Range("AY7").Value = 18
$ws.Range("AZ7").Value = 60
$ws.Range("BB7").Value = 175
$ws.Range("U8").Value = 1.67
$ws.Range("U9").Value = 1.73
$ws.Range("G10").Value = 1.91
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 2.63
$ws.Range("L10").Value = 4.5
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 8.5
$ws.Range("O10").Value = 1.36
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("X10").Value = 8.5
$ws.Range("Z10").Value = 17
$ws.Range("AO10").Value = 11
$ws.Range("AQ10").Value = 41
$ws.Range("M11").Value = 1.08
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("M12").Value = 1.08
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("M13").Value = 1.05
$ws.Range("O13").Value = 1.25
